$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 4386.86
$ws.Cells.Item(15, 9).Value = 4386.86
$ws.Cells.Item(15, 11).Value = 13160.58
$ws.Cells.Item(15, 13).Value = -12991.58
$ws.Cells.Item(17, 8).Value = 12502956
$ws.Cells.Item(17, 10).Value = 14288807
$ws.Cells.Item(17, 12).Value = 42866421
$ws.Cells.Item(17, 14).Value = -42866757
$ws.Cells.Item(19, 8).Value = 4830.8335
$ws.Cells.Item(19, 9).Value = 4158.4
$ws.Cells.Item(19, 10).Value = 5311.143
$ws.Cells.Item(19, 11).Value = 4158.4
$ws.Cells.Item(19, 12).Value = 5311.143
$ws.Cells.Item(19, 13).Value = -3983.4
$ws.Cells.Item(19, 14).Value = -5661.143
$ws.Cells.Item(86, 8).Value = 1807
$ws.Cells.Item(86, 9).Value = 1950
$ws.Cells.Item(86, 10).Value = 1449.5
$ws.Cells.Item(86, 11).Value = 1950
$ws.Cells.Item(86, 12).Value = 1449.5
$ws.Cells.Item(86, 13).Value = -827
$ws.Cells.Item(86, 14).Value = -3695.5
$ws.Cells.Item(89, 8).Value = 1807
$ws.Cells.Item(89, 9).Value = 1950
$ws.Cells.Item(89, 10).Value = 1449.5
$ws.Cells.Item(89, 11).Value = 9750
$ws.Cells.Item(89, 12).Value = 7247.5
$ws.Cells.Item(89, 13).Value = -4134
$ws.Cells.Item(89, 14).Value = -18479.5
$ws.Cells.Item(98, 8).Value = 1729.5
$ws.Cells.Item(98, 9).Value = 1748.0714
$ws.Cells.Item(98, 11).Value = 1748.0714
$ws.Cells.Item(98, 13).Value = -250.0714
$ws.Cells.Item(122, 8).Value = 1729.5
$ws.Cells.Item(122, 9).Value = 1748.0714
$ws.Cells.Item(122, 11).Value = 5244.2142
$ws.Cells.Item(122, 13).Value = -2794.2142
$ws.Cells.Item(137, 8).Value = 5222.4517
$ws.Cells.Item(137, 9).Value = 5292.875
$ws.Cells.Item(137, 10).Value = 4981
$ws.Cells.Item(137, 11).Value = 15878.625
$ws.Cells.Item(137, 12).Value = 14943
$ws.Cells.Item(137, 13).Value = -13328.625
$ws.Cells.Item(137, 14).Value = -20043
$ws.Cells.Item(138, 8).Value = 5266.5586
$ws.Cells.Item(138, 10).Value = 5017.7393
$ws.Cells.Item(138, 12).Value = 15053.2179
$ws.Cells.Item(138, 14).Value = -25333.2179

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10845.552
$ws.Cells.Item(32, 9).Value = 9123.483
$ws.Cells.Item(32, 11).Value = 9123.483
$ws.Cells.Item(32, 13).Value = -8836.483
$ws.Cells.Item(122, 8).Value = 2244
$ws.Cells.Item(122, 9).Value = 2096
$ws.Cells.Item(122, 11).Value = 6288
$ws.Cells.Item(122, 13).Value = -3838

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4201.5
$ws.Cells.Item(20, 9).Value = 4018.1
$ws.Cells.Item(20, 11).Value = 4018.1
$ws.Cells.Item(20, 13).Value = -3771.1
$ws.Cells.Item(22, 8).Value = 1139.6923
$ws.Cells.Item(22, 9).Value = 1026.3334
$ws.Cells.Item(22, 11).Value = 1026.3334
$ws.Cells.Item(22, 13).Value = -853.3334
$ws.Cells.Item(33, 8).Value = 1750
$ws.Cells.Item(33, 9).Value = 1750
$ws.Cells.Item(33, 11).Value = 1750
$ws.Cells.Item(33, 13).Value = -1414
$ws.Cells.Item(134, 8).Value = 1582.579
$ws.Cells.Item(134, 9).Value = 1582.579
$ws.Cells.Item(134, 11).Value = 4747.737
$ws.Cells.Item(134, 13).Value = -2212.737

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2845.6667
$ws.Cells.Item(16, 9).Value = 2740.7273
$ws.Cells.Item(16, 10).Value = 4000
$ws.Cells.Item(16, 11).Value = 2740.7273
$ws.Cells.Item(16, 12).Value = 4000
$ws.Cells.Item(16, 13).Value = -2453.7273
$ws.Cells.Item(16, 14).Value = -4574
$ws.Cells.Item(31, 8).Value = 4411.7646
$ws.Cells.Item(31, 9).Value = 2468.923
$ws.Cells.Item(31, 10).Value = 10726
$ws.Cells.Item(31, 11).Value = 2468.923
$ws.Cells.Item(31, 12).Value = 10726
$ws.Cells.Item(31, 13).Value = -2173.923
$ws.Cells.Item(31, 14).Value = -11316
$ws.Cells.Item(34, 8).Value = 4411.7646
$ws.Cells.Item(34, 9).Value = 2468.923
$ws.Cells.Item(34, 10).Value = 10726
$ws.Cells.Item(34, 11).Value = 2468.923
$ws.Cells.Item(34, 12).Value = 10726
$ws.Cells.Item(34, 13).Value = -2266.923
$ws.Cells.Item(34, 14).Value = -11130
$ws.Cells.Item(62, 8).Value = 6201.25
$ws.Cells.Item(62, 9).Value = 6466.3335
$ws.Cells.Item(62, 10).Value = 5406
$ws.Cells.Item(62, 11).Value = 6466.3335
$ws.Cells.Item(62, 12).Value = 5406
$ws.Cells.Item(62, 13).Value = -5842.3335
$ws.Cells.Item(62, 14).Value = -6654
$ws.Cells.Item(65, 8).Value = 6201.25
$ws.Cells.Item(65, 9).Value = 6466.3335
$ws.Cells.Item(65, 10).Value = 5406
$ws.Cells.Item(65, 11).Value = 32331.6675
$ws.Cells.Item(65, 12).Value = 27030
$ws.Cells.Item(65, 13).Value = -29211.6675
$ws.Cells.Item(65, 14).Value = -33270
$ws.Cells.Item(113, 8).Value = 2845.6667
$ws.Cells.Item(113, 9).Value = 2740.7273
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 2740.7273
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = -570.7273
$ws.Cells.Item(113, 14).Value = -8340
$ws.Cells.Item(132, 8).Value = 3100.3333
$ws.Cells.Item(132, 10).Value = 7885.7144
$ws.Cells.Item(132, 12).Value = 23657.1432
$ws.Cells.Item(132, 14).Value = -28717.1432
$ws.Cells.Item(134, 8).Value = 1980.4562
$ws.Cells.Item(134, 9).Value = 1125.2354
$ws.Cells.Item(134, 11).Value = 3375.7062
$ws.Cells.Item(134, 13).Value = -840.7062000000001
$ws.Cells.Item(141, 8).Value = 251297.5
$ws.Cells.Item(141, 9).Value = 47647.75
$ws.Cells.Item(141, 10).Value = 319180.75
$ws.Cells.Item(141, 11).Value = 47647.75
$ws.Cells.Item(141, 12).Value = 319180.75
$ws.Cells.Item(141, 13).Value = -42467.75
$ws.Cells.Item(141, 14).Value = -329540.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 29280654
$ws.Cells.Item(4, 9).Value = 28164004
$ws.Cells.Item(4, 11).Value = 84492012
$ws.Cells.Item(4, 13).Value = -84491900
$ws.Cells.Item(10, 8).Value = 35.25
$ws.Cells.Item(10, 9).Value = 40.333332
$ws.Cells.Item(10, 10).Value = 20
$ws.Cells.Item(10, 11).Value = 120.999996
$ws.Cells.Item(10, 12).Value = 60
$ws.Cells.Item(10, 13).Value = 18.000004
$ws.Cells.Item(10, 14).Value = -338
$ws.Cells.Item(86, 8).Value = 424.5
$ws.Cells.Item(86, 9).Value = 424.5
$ws.Cells.Item(86, 11).Value = 1273.5
$ws.Cells.Item(86, 13).Value = -87.5
$ws.Cells.Item(89, 8).Value = 424.5
$ws.Cells.Item(89, 9).Value = 424.5
$ws.Cells.Item(89, 11).Value = 3820.5
$ws.Cells.Item(89, 13).Value = 2107.5
$ws.Cells.Item(113, 8).Value = 6135.4546
$ws.Cells.Item(113, 10).Value = 7310
$ws.Cells.Item(113, 12).Value = 21930
$ws.Cells.Item(113, 14).Value = -26270
$ws.Cells.Item(122, 8).Value = 4458.1113
$ws.Cells.Item(122, 10).Value = 4803.222
$ws.Cells.Item(122, 12).Value = 43228.998
$ws.Cells.Item(122, 14).Value = -48128.998
$ws.Cells.Item(128, 8).Value = 209993.67
$ws.Cells.Item(128, 9).Value = 209993.67
$ws.Cells.Item(128, 11).Value = 629981.01
$ws.Cells.Item(128, 13).Value = -625001.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 26269.666
$ws.Cells.Item(22, 10).Value = 36904.5
$ws.Cells.Item(22, 12).Value = 36904.5
$ws.Cells.Item(22, 14).Value = -37962.5
$ws.Cells.Item(70, 8).Value = 8969.700000000001
$ws.Cells.Item(70, 9).Value = 3996.3333
$ws.Cells.Item(70, 11).Value = 3996.3333
$ws.Cells.Item(70, 13).Value = -3726.3333
$ws.Cells.Item(73, 8).Value = 8969.700000000001
$ws.Cells.Item(73, 9).Value = 3996.3333
$ws.Cells.Item(73, 11).Value = 3996.3333
$ws.Cells.Item(73, 13).Value = -3060.3333
$ws.Cells.Item(80, 8).Value = 4223.75
$ws.Cells.Item(80, 9).Value = 1944.5
$ws.Cells.Item(80, 11).Value = 1944.5
$ws.Cells.Item(80, 13).Value = -946.5
$ws.Cells.Item(83, 8).Value = 4223.75
$ws.Cells.Item(83, 9).Value = 1944.5
$ws.Cells.Item(83, 11).Value = 9722.5
$ws.Cells.Item(83, 13).Value = -4730.5
$ws.Cells.Item(132, 8).Value = 3330.5417
$ws.Cells.Item(132, 9).Value = 2706.75
$ws.Cells.Item(132, 11).Value = 8120.25
$ws.Cells.Item(132, 13).Value = -5590.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 60516
$ws.Cells.Item(82, 9).Value = 1970.4546
$ws.Cells.Item(82, 10).Value = 167849.5
$ws.Cells.Item(82, 11).Value = 1970.4546
$ws.Cells.Item(82, 12).Value = 167849.5
$ws.Cells.Item(82, 13).Value = -1609.4546
$ws.Cells.Item(82, 14).Value = -168571.5
$ws.Cells.Item(85, 8).Value = 60516
$ws.Cells.Item(85, 9).Value = 1970.4546
$ws.Cells.Item(85, 10).Value = 167849.5
$ws.Cells.Item(85, 11).Value = 1970.4546
$ws.Cells.Item(85, 12).Value = 167849.5
$ws.Cells.Item(85, 13).Value = -722.4546
$ws.Cells.Item(85, 14).Value = -170345.5
$ws.Cells.Item(132, 8).Value = 9843.627
$ws.Cells.Item(132, 9).Value = 9657.361000000001
$ws.Cells.Item(132, 11).Value = 28972.083
$ws.Cells.Item(132, 13).Value = -26442.083
$ws.Cells.Item(135, 8).Value = 66431.664
$ws.Cells.Item(135, 10).Value = 66431.664
$ws.Cells.Item(135, 12).Value = 66431.664
$ws.Cells.Item(135, 14).Value = -76571.664
$ws.Cells.Item(136, 8).Value = 83341064
$ws.Cells.Item(136, 9).Value = 45462524
$ws.Cells.Item(136, 11).Value = 136387572
$ws.Cells.Item(136, 13).Value = -136385022
$ws.Cells.Item(140, 8).Value = 74429
$ws.Cells.Item(140, 10).Value = 74429
$ws.Cells.Item(140, 12).Value = 74429
$ws.Cells.Item(140, 14).Value = -84789

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 461655.66
$ws.Cells.Item(5, 10).Value = 461655.66
$ws.Cells.Item(5, 12).Value = 461655.66
$ws.Cells.Item(5, 14).Value = -461879.66
$ws.Cells.Item(13, 8).Value = 995
$ws.Cells.Item(13, 9).Value = 995
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 995
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = -855
$ws.Cells.Item(13, 14).Value = $null
$ws.Cells.Item(41, 8).Value = 19442.666
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 19442.666
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 19442.666
$ws.Cells.Item(41, 13).Value = $null
$ws.Cells.Item(41, 14).Value = -20222.666
$ws.Cells.Item(46, 8).Value = 81476.08
$ws.Cells.Item(46, 10).Value = 83275.164
$ws.Cells.Item(46, 12).Value = 83275.164
$ws.Cells.Item(46, 14).Value = -83737.164
$ws.Cells.Item(107, 8).Value = 1721
$ws.Cells.Item(107, 9).Value = 1090.3334
$ws.Cells.Item(107, 10).Value = 2581
$ws.Cells.Item(107, 11).Value = 3271.0002
$ws.Cells.Item(107, 12).Value = 7743
$ws.Cells.Item(107, 13).Value = -1351.0002
$ws.Cells.Item(107, 14).Value = -11583
$ws.Cells.Item(134, 8).Value = 81476.08
$ws.Cells.Item(134, 10).Value = 83275.164
$ws.Cells.Item(134, 12).Value = 249825.492
$ws.Cells.Item(134, 14).Value = -254895.492
